$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking Price text (e.g. "2.960") keeps its exact
# text representation instead of being auto-converted to a number.
$textFormatCells = @("D5","D6","D7","D8","D9","D10","D11","D13","D14","D16","D18","D20","D22","D23","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '27.019.95'
$ws.Range("E2").Value = '  -2.45%  '

$ws.Range("D3").Value = '1.818.83'
$ws.Range("E3").Value = '  -1.42%  '

$ws.Range("E4").Value = '  -1.27%  '

$ws.Range("D5").Value = '310.64'
$ws.Range("E5").Value = '  -2.64%  '

$ws.Range("D6").Value = '0.9999'
$ws.Range("E6").Value = '  -1.10%  '

$ws.Range("D7").Value = '0.4219'
$ws.Range("E7").Value = '  -2.08%  '

$ws.Range("D8").Value = '0.3671'
$ws.Range("E8").Value = '  -2.03%  '

$ws.Range("D9").Value = '0.07208'
$ws.Range("E9").Value = '  -1.90%  '

$ws.Range("D10").Value = '0.8390'
$ws.Range("E10").Value = '  -4.30%  '

$ws.Range("D11").Value = '20.76'
$ws.Range("E11").Value = '  -3.88%  '

$ws.Range("D12").Value = '1.813.59'
$ws.Range("E12").Value = '  -1.79%  '

$ws.Range("D13").Value = '6.644'
$ws.Range("E13").Value = '  -1.10%  '

$ws.Range("D14").Value = '0.07073'
$ws.Range("E14").Value = '  -0.86%  '

$ws.Range("E15").Value = '  -3.05%  '

$ws.Range("D16").Value = '89.32'
$ws.Range("E16").Value = '  +0.35%  '

$ws.Range("E17").Value = '  -1.41%  '

$ws.Range("D18").Value = '0.000008796'
$ws.Range("E18").Value = '  -2.12%  '

$ws.Range("E19").Value = '  -0.98%  '

$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '14.92'
$ws.Range("E20").Value = '  -3.55%  '

$ws.Range("B21").Value = 'WrappedBTC'
$ws.Range("C21").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D21").Value = '27.061.45'
$ws.Range("E21").Value = '  -2.28%  '

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '5.119'
$ws.Range("E22").Value = '  -1.87%  '

$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").Value = '10.84'
$ws.Range("E23").Value = '  -2.18%  '

$ws.Range("B24").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D24").Value = '2.039.71'
$ws.Range("E24").Value = '  -1.83%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '1.979'
$ws.Range("E25").Value = '  -1.33%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '151.81'
$ws.Range("E26").Value = '  -2.16%  '

$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = '2.216'
$ws.Range("E27").Value = '  +1.41%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '18.23'
$ws.Range("E28").Value = '  -2.20%  '

$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").Value = '5.214'
$ws.Range("E29").Value = '  -3.11%  '

$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").Value = '116.24'
$ws.Range("E30").Value = '  -2.48%  '

$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = '0.08744'
$ws.Range("E31").Value = '  -2.16%  '

$ws.Range("B32").Value = 'ARBITRUM'
$ws.Range("C32").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D32").Value = '1.177'
$ws.Range("E32").Value = '  -4.32%  '

$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").Value = '2.960'
$ws.Range("E33").Value = '  +1.33%  '

$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = '0.7374'
$ws.Range("E34").Value = '  -5.15%  '

$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").Value = '4.407'
$ws.Range("E35").Value = '  -3.24%  '

$ws.Range("B36").Value = 'Frax'
$ws.Range("C36").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D36").Value = '0.9995'
$ws.Range("E36").Value = '  -1.31%  '

$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '1.090'
$ws.Range("E37").Value = '  -3.94%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.01952'
$ws.Range("E38").Value = '  -1.07%  '

$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '0.05231'
$ws.Range("E39").Value = '  -1.94%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '7.269'
$ws.Range("E40").Value = '  -0.58%  '

$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = '2.867'
$ws.Range("E41").Value = '  -1.12%  '

$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '0.1687'
$ws.Range("E42").Value = '  +0.02%  '

$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '0.5026'
$ws.Range("E43").Value = '  -1.77%  '

$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").Value = '8.581'
$ws.Range("E44").Value = '  -2.49%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '10.47'
$ws.Range("E45").Value = '  -1.87%  '

$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").Value = '106.11'
$ws.Range("E46").Value = '  -2.59%  '

$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '0.4708'
$ws.Range("E47").Value = '  -1.00%  '

$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").Value = '0.9991'
$ws.Range("E48").Value = '  -1.30%  '

$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '0.06344'
$ws.Range("E49").Value = '  -2.07%  '

$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '1.644'
$ws.Range("E50").Value = '  -2.71%  '

$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").Value = '1.875'
$ws.Range("E51").Value = '  +1.62%  '
